$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F13").Value = '135_product_information'
$ws.Range("F17").Value = '135_product_information'
$ws.Range("F54").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F59").Value = 'ppe'
$ws.Range("F68").Value = 'pollinator'
$ws.Range("F69").Value = 'pollinator'
$ws.Range("F76").Value = 'pollinator'
$ws.Range("F77").Value = 'pollinator'
$ws.Range("F78").Value = 'pollinator'
$ws.Range("F79").Value = 'pollinator'
$ws.Range("F81").Value = 'pollinator'
$ws.Range("F82").Value = 'pollinator'
$ws.Range("F84").Value = 'pollinator'
$ws.Range("F85").Value = 'pollinator'
$ws.Range("F86").Value = 'pollinator'
$ws.Range("F87").Value = 'pollinator'
$ws.Range("F88").Value = 'pollinator'
$ws.Range("F93").Value = 'pollinator'
$ws.Range("F95").Value = 'pollinator'
$ws.Range("F96").Value = 'pollinator'
$ws.Range("F98").Value = 'application instructions'
$ws.Range("F123").Value = '135_product_information'
$ws.Range("F124").Value = 'use restrictions'
$ws.Range("F137").Value = 'mixing'
$ws.Range("F138").Value = 'off target movement'
$ws.Range("F140").Value = 'application instructions'
$ws.Range("F141").Value = 'irrigation || application instructions || chemigation'
$ws.Range("F142").Value = 'irrigation'
$ws.Range("F143").Value = 'irrigation'
$ws.Range("F144").Value = 'safety procedures || irrigation'
$ws.Range("F145").Value = 'application instructions'
$ws.Range("F146").Value = 'irrigation'
$ws.Range("F147").Value = 'irrigation'
$ws.Range("F148").Value = 'irrigation'
$ws.Range("F214").Value = 'application instructions'
$ws.Range("F326").Value = 'use restrictions'
$ws.Range("F328").Value = '154_pesticide_storage'
$ws.Range("F343").Value = '135_product_information'
$ws.Range("F351").Value = '135_product_information'
$ws.Range("F355").Value = '135_product_information'
$ws.Range("F366").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F374").Value = 'ppe'
$ws.Range("F383").Value = 'pollinator'
$ws.Range("F384").Value = 'pollinator'
$ws.Range("F399").Value = 'pollinator'
$ws.Range("F401").Value = 'pollinator'
$ws.Range("F402").Value = 'pollinator'
$ws.Range("F404").Value = 'pollinator'
$ws.Range("F406").Value = 'application instructions'
$ws.Range("F407").Value = '134_non-agriculture_use_requirements'
$ws.Range("F431").Value = '135_product_information'
$ws.Range("F432").Value = 'use restrictions'
$ws.Range("F434").Value = 'mixing'
$ws.Range("F435").Value = 'off target movement'
$ws.Range("F437").Value = 'application instructions'
$ws.Range("F438").Value = 'irrigation || application instructions || chemigation'
$ws.Range("F439").Value = 'irrigation'
$ws.Range("F440").Value = 'irrigation'
$ws.Range("F441").Value = 'safety procedures || irrigation'
$ws.Range("F442").Value = 'application instructions'
$ws.Range("F443").Value = 'irrigation'
$ws.Range("F444").Value = 'irrigation'
$ws.Range("F445").Value = 'irrigation'
$ws.Range("F446").Value = 'application instructions'
$ws.Range("F447").Value = 'use restrictions'
$ws.Range("F466").Value = 'use restrictions'
$ws.Range("F472").Value = 'use restrictions'
$ws.Range("F478").Value = 'use restrictions'
$ws.Range("F501").Value = 'application instructions'
$ws.Range("F502").Value = 'use restrictions'
$ws.Range("F512").Value = 'use restrictions'
$ws.Range("F522").Value = 'use restrictions'
$ws.Range("F557").Value = 'use restrictions'
$ws.Range("F560").Value = 'application instructions'
$ws.Range("F574").Value = 'application instructions'
$ws.Range("F584").Value = 'application instructions'
$ws.Range("F598").Value = 'use restrictions'
$ws.Range("F600").Value = '154_pesticide_storage'
